$d = $word.ActiveDocument

# 1. Replace the title text "2.2 - Debate I" with "Placeholder - Check Back Later".
#    Using Range.Text (rather than Find.Execute's ReplaceWith) keeps the
#    xml:space="preserve" attribute on the surviving run's <w:t>.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Text = "Placeholder - Check Back Later"

# 2. Remove the trailing " " and ":::" runs that follow the
#    "...general edification later." sentence inside the table cell.
$d.Content.Find.Execute(" :::", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
